$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stale "_GoBack" bookmark that currently sits right after
#    the word "result" in the result_queue_count_threshold_severity
#    paragraph (in the Constants table). Word re-homes this bookmark to
#    the location of the most recent edit whenever the document is
#    saved, so it must be removed from its old spot before the new edit
#    is made. The identical wording also appears earlier in the
#    Generics table, so Find alone cannot disambiguate the two -- the
#    bookmark's exact character offset (in the freshly-loaded document)
#    is used instead.
# ---------------------------------------------------------------------
$goBackAnchorText = "An alert with severity 'result_queue_count_threshold_severity' will be issued if result"
$anchorProbe = $d.Content
$anchorFound = $false
while ($anchorProbe.Find.Execute($goBackAnchorText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)) {
    $precedingRange = $d.Range([Math]::Max(0, $anchorProbe.Start - 400), $anchorProbe.Start)
    if ($precedingRange.Text -match "(?<![A-Z])C_RESULT_QUEUE_COUNT_THRESHOLD(?!_)") {
        $anchorFound = $true
        break
    }
    $anchorProbe.Collapse(0)
}
if ($anchorFound) {
    $wordEnd = $anchorProbe.End
    # Span a small range that straddles the collapsed bookmark location
    # (right after "result"), delete it and retype the same text -- this
    # removes any bookmark anchored inside the deleted span.
    $spanStart = $wordEnd - 2
    $spanEnd = $wordEnd + 10
    $span = $d.Range($spanStart, $spanEnd)
    $savedText = $span.Text
    $span.Delete()
    $reinsertPoint = $d.Range($spanStart, $spanStart)
    $reinsertPoint.InsertAfter($savedText)
}

# ---------------------------------------------------------------------
# 2) Fix the VVC QuickRef wording: replace the old two-package mention
#    with the single, correctly named package.
# ---------------------------------------------------------------------
$editRange = $d.Content
$oldText = "uvvm_vvc_framework.uvvm_methods_pkg and uvvm_vvc_framework.uvvm_support_pkg"
$newText = "uvvm_vvc_framework.td_vvc_framework_common_methods_pkg"
$editRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Word marks the spot of this latest edit with a new "_GoBack"
#    bookmark, collapsed right after "uvvm_vvc_framework.".
# ---------------------------------------------------------------------
$locateRange = $d.Content
$found = $locateRange.Find.Execute("uvvm_vvc_framework.td_vvc_framework_common_methods_pkg", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $bookmarkPos = $locateRange.Start + [int]("uvvm_vvc_framework.".Length)
    $bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
}

# ---------------------------------------------------------------------
# 4) Update the cached result of the footer's DATE field.
# ---------------------------------------------------------------------
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections($s)
    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        $footer = $section.Footers($f)
        if ($footer.Exists) {
            $footerRange = $footer.Range
            $footerRange.Find.Execute("2018-07-31", $true, $false, $false, $false, $false, $true, 1, $false, "2018-11-19", 2) | Out-Null
        }
    }
}
